$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.91845900000001
$ws.Range("N2").Value = 131.755377
$ws.Range("O2").Value = 0.8150909120558799
$ws.Range("P2").Value = 0.81509091205588
$ws.Range("Q2").Value = 959.6081693464847
$ws.Range("R2").Value = 8636.473524118363
$ws.Range("S2").Value = 0.04091771776304753
$ws.Range("T2").Value = 0.04091771776304753

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.100310333333333
$ws.Range("N3").Value = 6.300930999999999
$ws.Range("O3").Value = 0.03898005312975703
$ws.Range("P3").Value = 0.03898005312975703
$ws.Range("Q3").Value = 45.89129491154289
$ws.Range("R3").Value = 413.021654203886
$ws.Range("S3").Value = 0.001956806030788685
$ws.Range("T3").Value = 0.001956806030788685

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.8629
$ws.Range("N4").Value = 23.5887
$ws.Range("O4").Value = 0.1459290348143631
$ws.Range("P4").Value = 0.1459290348143631
$ws.Range("Q4").Value = 171.8025460491333
$ws.Range("R4").Value = 1546.2229144422
$ws.Range("S4").Value = 0.007325665114959212
$ws.Range("T4").Value = 0.007325665114959212

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 43.91845900000001
$ws.Range("N5").Value = 131.755377
$ws.Range("O5").Value = 0.8150909120558799
$ws.Range("P5").Value = 0.81509091205588
$ws.Range("Q5").Value = 16910.91113118322
$ws.Range("R5").Value = 152198.200180649
$ws.Range("S5").Value = 0.721081698640573
$ws.Range("T5").Value = 0.7210816986405731

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.100310333333333
$ws.Range("N6").Value = 6.300930999999999
$ws.Range("O6").Value = 0.03898005312975703
$ws.Range("P6").Value = 0.03898005312975703
$ws.Range("Q6").Value = 808.7296823166268
$ws.Range("R6").Value = 7278.567140849641
$ws.Range("S6").Value = 0.03448425507899418
$ws.Range("T6").Value = 0.03448425507899418

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.8629
$ws.Range("N7").Value = 23.5887
$ws.Range("O7").Value = 0.1459290348143631
$ws.Range("P7").Value = 0.1459290348143631
$ws.Range("Q7").Value = 3027.629068983967
$ws.Range("R7").Value = 27248.6616208557
$ws.Range("S7").Value = 0.129098183709974
$ws.Range("T7").Value = 0.129098183709974

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 43.91845900000001
$ws.Range("N8").Value = 131.755377
$ws.Range("O8").Value = 0.8150909120558799
$ws.Range("P8").Value = 0.81509091205588
$ws.Range("Q8").Value = 1245.109349591861
$ws.Range("R8").Value = 11205.98414632675
$ws.Range("S8").Value = 0.05309149565225935
$ws.Range("T8").Value = 0.05309149565225936

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.100310333333333
$ws.Range("N9").Value = 6.300930999999999
$ws.Range("O9").Value = 0.03898005312975703
$ws.Range("P9").Value = 0.03898005312975703
$ws.Range("Q9").Value = 59.54480399864966
$ws.Range("R9").Value = 535.9032359878469
$ws.Range("S9").Value = 0.002538992019974153
$ws.Range("T9").Value = 0.002538992019974153

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.8629
$ws.Range("N10").Value = 23.5887
$ws.Range("O10").Value = 0.1459290348143631
$ws.Range("P10").Value = 0.1459290348143631
$ws.Range("Q10").Value = 222.9169813291
$ws.Range("R10").Value = 2006.2528319619
$ws.Range("S10").Value = 0.009505185989429867
$ws.Range("T10").Value = 0.009505185989429867
